$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for "「信者」..." (row 644) was removed from the post list.
# Deleting the entire row shifts all subsequent rows up by one.
$ws.Rows.Item(644).Delete()
